$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3547.4167
$ws.Range("I15").Value = 3547.4167
$ws.Range("K15").Value = 10642.2501
$ws.Range("M15").Value = -10473.2501
$ws.Range("H17").Value = 1318.919
$ws.Range("J17").Value = 1323.1428
$ws.Range("L17").Value = 3969.4284
$ws.Range("N17").Value = -4305.428400000001
$ws.Range("H19").Value = 948
$ws.Range("J19").Value = 1051
$ws.Range("L19").Value = 1051
$ws.Range("N19").Value = -1401
$ws.Range("H31").Value = 33
$ws.Range("I31").Value = 33
$ws.Range("K31").Value = 99
$ws.Range("M31").Value = 131
$ws.Range("H41").Value = 162.4
$ws.Range("I41").Value = 125.85714
$ws.Range("J41").Value = 247.66667
$ws.Range("K41").Value = 125.85714
$ws.Range("L41").Value = 247.66667
$ws.Range("M41").Value = 314.14286
$ws.Range("N41").Value = -1127.66667
$ws.Range("H43").Value = 2095
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 2399
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 2399
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -2537
$ws.Range("H96").Value = 695
$ws.Range("I96").Value = 695
$ws.Range("K96").Value = 2085
$ws.Range("M96").Value = -712
$ws.Range("H103").Value = 3157
$ws.Range("I103").Value = 3459.6667
$ws.Range("J103").Value = 2551.6667
$ws.Range("K103").Value = 10379.0001
$ws.Range("L103").Value = 7655.000100000001
$ws.Range("M103").Value = -9793.000100000001
$ws.Range("N103").Value = -8827.000100000001
$ws.Range("H112").Value = 1837.05
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1837.05
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 5511.15
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -7727.15
$ws.Range("H129").Value = 2273.75
$ws.Range("J129").Value = 3660.75
$ws.Range("L129").Value = 10982.25
$ws.Range("N129").Value = -20982.25
$ws.Range("H137").Value = 2493.121
$ws.Range("I137").Value = 1522.5555
$ws.Range("J137").Value = 6860.6665
$ws.Range("K137").Value = 4567.666499999999
$ws.Range("L137").Value = 20581.9995
$ws.Range("M137").Value = -2017.666499999999
$ws.Range("N137").Value = -25681.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1576.5714
$ws.Range("I45").Value = 1341.7778
$ws.Range("K45").Value = 1341.7778
$ws.Range("M45").Value = -964.7778000000001
$ws.Range("H61").Value = 1663.7667
$ws.Range("I61").Value = 1089.2727
$ws.Range("K61").Value = 1089.2727
$ws.Range("M61").Value = -877.2727
$ws.Range("H74").Value = 2591.673
$ws.Range("I74").Value = 2117.4055
$ws.Range("K74").Value = 2117.4055
$ws.Range("M74").Value = -1243.4055
$ws.Range("H77").Value = 2591.673
$ws.Range("I77").Value = 2117.4055
$ws.Range("K77").Value = 10587.0275
$ws.Range("M77").Value = -6219.0275
$ws.Range("H88").Value = 4259.8184
$ws.Range("I88").Value = 3153.3333
$ws.Range("J88").Value = 5587.6
$ws.Range("K88").Value = 3153.3333
$ws.Range("L88").Value = 5587.6
$ws.Range("M88").Value = -2747.3333
$ws.Range("N88").Value = -6399.6
$ws.Range("H91").Value = 4259.8184
$ws.Range("I91").Value = 3153.3333
$ws.Range("J91").Value = 5587.6
$ws.Range("K91").Value = 3153.3333
$ws.Range("L91").Value = 5587.6
$ws.Range("M91").Value = -1749.3333
$ws.Range("N91").Value = -8395.6
$ws.Range("H122").Value = 4035.8125
$ws.Range("I122").Value = 4068.8572
$ws.Range("J122").Value = 3804.5
$ws.Range("K122").Value = 12206.5716
$ws.Range("L122").Value = 11413.5
$ws.Range("M122").Value = -9756.571599999999
$ws.Range("N122").Value = -16313.5
$ws.Range("H132").Value = 2037.8334
$ws.Range("I132").Value = 1415.2727
$ws.Range("J132").Value = 3749.875
$ws.Range("K132").Value = 4245.8181
$ws.Range("L132").Value = 11249.625
$ws.Range("M132").Value = -1715.8181
$ws.Range("N132").Value = -16309.625
$ws.Range("H136").Value = 1663.7667
$ws.Range("I136").Value = 1089.2727
$ws.Range("K136").Value = 3267.8181
$ws.Range("M136").Value = -717.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 47619950
$ws.Range("I94").Value = 71429360
$ws.Range("J94").Value = 1145.5
$ws.Range("K94").Value = 71429360
$ws.Range("L94").Value = 1145.5
$ws.Range("M94").Value = -71428909
$ws.Range("N94").Value = -2047.5
$ws.Range("H107").Value = 6994461.5
$ws.Range("I107").Value = 12821634
$ws.Range("J107").Value = 1855
$ws.Range("K107").Value = 12821634
$ws.Range("L107").Value = 1855
$ws.Range("M107").Value = -12819714
$ws.Range("N107").Value = -5695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 265.8
$ws.Range("I7").Value = 174.75
$ws.Range("J7").Value = 402.375
$ws.Range("K7").Value = 174.75
$ws.Range("L7").Value = 402.375
$ws.Range("M7").Value = -61.75
$ws.Range("N7").Value = -628.375
$ws.Range("H31").Value = 3546.8096
$ws.Range("I31").Value = 2403.55
$ws.Range("J31").Value = 4586.136
$ws.Range("K31").Value = 2403.55
$ws.Range("L31").Value = 4586.136
$ws.Range("M31").Value = -2108.55
$ws.Range("N31").Value = -5176.136
$ws.Range("H34").Value = 3546.8096
$ws.Range("I34").Value = 2403.55
$ws.Range("J34").Value = 4586.136
$ws.Range("K34").Value = 2403.55
$ws.Range("L34").Value = 4586.136
$ws.Range("M34").Value = -2201.55
$ws.Range("N34").Value = -4990.136
$ws.Range("H94").Value = 1358
$ws.Range("I94").Value = 888.6
$ws.Range("J94").Value = 1651.375
$ws.Range("K94").Value = 888.6
$ws.Range("L94").Value = 1651.375
$ws.Range("M94").Value = -437.6
$ws.Range("N94").Value = -2553.375
$ws.Range("H107").Value = 2175211.5
$ws.Range("I107").Value = 3334406.2
$ws.Range("J107").Value = 1721.375
$ws.Range("K107").Value = 3334406.2
$ws.Range("L107").Value = 1721.375
$ws.Range("M107").Value = -3332486.2
$ws.Range("N107").Value = -5561.375
$ws.Range("H132").Value = 12351718
$ws.Range("I132").Value = 3160.7896
$ws.Range("K132").Value = 9482.3688
$ws.Range("M132").Value = -6952.3688
$ws.Range("H134").Value = 2053.0466
$ws.Range("I134").Value = 1918.8572
$ws.Range("J134").Value = 2640.125
$ws.Range("K134").Value = 5756.571599999999
$ws.Range("L134").Value = 7920.375
$ws.Range("M134").Value = -3221.571599999999
$ws.Range("N134").Value = -12990.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 410.125
$ws.Range("I34").Value = 325.85715
$ws.Range("K34").Value = 977.5714499999999
$ws.Range("M34").Value = -893.5714499999999
$ws.Range("H56").Value = 7296.375
$ws.Range("I56").Value = 7296.375
$ws.Range("K56").Value = 7296.375
$ws.Range("M56").Value = -6766.375
$ws.Range("H92").Value = 899.6
$ws.Range("J92").Value = 899.6
$ws.Range("L92").Value = 2698.8
$ws.Range("N92").Value = -5194.8
$ws.Range("H107").Value = 3000
$ws.Range("J107").Value = 3150.4443
$ws.Range("L107").Value = 9451.332900000001
$ws.Range("N107").Value = -13291.3329
$ws.Range("H131").Value = 2357.258
$ws.Range("I131").Value = 6015
$ws.Range("K131").Value = 18045
$ws.Range("M131").Value = -13005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1501.2778
$ws.Range("I97").Value = 1289.3636
$ws.Range("K97").Value = 1289.3636
$ws.Range("M97").Value = -793.3635999999999
$ws.Range("H107").Value = 7403
$ws.Range("I107").Value = 439.2
$ws.Range("K107").Value = 439.2
$ws.Range("M107").Value = 1480.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4029.1667
$ws.Range("I22").Value = 4383.4443
$ws.Range("J22").Value = 2966.3333
$ws.Range("K22").Value = 4383.4443
$ws.Range("L22").Value = 2966.3333
$ws.Range("M22").Value = -4088.4443
$ws.Range("N22").Value = -3556.3333
$ws.Range("H27").Value = 4029.1667
$ws.Range("I27").Value = 4383.4443
$ws.Range("J27").Value = 2966.3333
$ws.Range("K27").Value = 4383.4443
$ws.Range("L27").Value = 2966.3333
$ws.Range("M27").Value = -4276.4443
$ws.Range("N27").Value = -3180.3333
$ws.Range("H122").Value = 16286.529
$ws.Range("I122").Value = 14189.4
$ws.Range("K122").Value = 42568.2
$ws.Range("M122").Value = -40118.2
$ws.Range("H132").Value = 7643.5713
$ws.Range("I132").Value = 4334.9443
$ws.Range("K132").Value = 13004.8329
$ws.Range("M132").Value = -10474.8329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4529.6665
$ws.Range("I62").Value = 3882.8333
$ws.Range("J62").Value = 5823.3335
$ws.Range("K62").Value = 3882.8333
$ws.Range("L62").Value = 5823.3335
$ws.Range("M62").Value = -3258.8333
$ws.Range("N62").Value = -7071.3335
$ws.Range("H65").Value = 4529.6665
$ws.Range("I65").Value = 3882.8333
$ws.Range("J65").Value = 5823.3335
$ws.Range("K65").Value = 19414.1665
$ws.Range("L65").Value = 29116.6675
$ws.Range("M65").Value = -16294.1665
$ws.Range("N65").Value = -35356.6675
$ws.Range("H81").Value = 6459
$ws.Range("I81").Value = 7147.5
$ws.Range("K81").Value = 14295
$ws.Range("M81").Value = -13234
$ws.Range("H84").Value = 6459
$ws.Range("I84").Value = 7147.5
$ws.Range("K84").Value = 71475
$ws.Range("M84").Value = -66171
$ws.Range("H100").Value = 166671980
$ws.Range("I100").Value = 5979
$ws.Range("K100").Value = 11958
$ws.Range("M100").Value = -11417
